# Add a new "Trim Whitespace" worksheet to the GOMS workbook, modeled on the
# existing "Special Character Removal" sheet, then populate it with the
# Trim-Whitespace task rows (per the author's Python IDE notes).

$wb = $excel.ActiveWorkbook

# Start from a duplicate of "Special Character Removal" so the new sheet
# inherits identical formatting (fonts/styles/row heights) placed right
# after it, then rename + retarget the content.
$src = $wb.Worksheets.Item("Special Character Removal")
$src.Copy($null, $src)

$ws = $wb.Worksheets.Item($src.Index + 1)
$ws.Name = "Trim Whitespace"

# The source sheet has an extra "Remove Characters" row (row 4) that Trim
# Whitespace doesn't need - drop it so we end up with 5 rows total.
$ws.Rows.Item(4).Delete()

# Re-point the remaining rows at the Trim Whitespace content.
$ws.Range("A3").Value = "Trim Whitespace"
$ws.Range("C3").Value = "df = df.applymap(lambda x: x.strip() if isinstance(x, str) else x)"

$ws.Range("C4").Value = "df.head() to check trimmed strings"

$ws.Range("B5").Value = "8 min"

# Match the saved selection on the new sheet.
$ws.Range("L11").Select()

# Make the new sheet the active one, matching the author's final view state.
$ws.Activate()
